$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing Wednesday hours for the week on row 8
$ws.Range("G8").Value = 6.5

# Move the active selection to the cell that was just edited
$ws.Range("G8").Select()
